$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition listings)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 8102
$wsExhibit.Range("F5").Value = 5910
$wsExhibit.Range("F7").Value = 94
$wsExhibit.Range("F11").Value = 437

# Sheet "全部类型" (all types, combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 8102
$wsAll.Range("F5").Value = 5910
$wsAll.Range("F7").Value = 94
$wsAll.Range("F15").Value = 437
